$wb = $excel.ActiveWorkbook

# ============================================================
# Sheet 1: Active Signals
# ============================================================
$ws1 = $wb.Worksheets.Item("Active Signals")

# Insert a new row at position 5 (old row 5 USDCHF/BUY shifts down to row 6)
$ws1.Rows.Item(5).Insert()

# The new row starts out unformatted; give it the same look as the other
# data rows (bordered cells, red Signal-column fill) before filling it in.
$ws1.Range("A2:J2").Copy()
$ws1.Range("A5:J5").PasteSpecial(-4122)
$ws1.Range("C2").Copy()
$ws1.Range("C5").PasteSpecial(-4122)

# Row 2: NZDUSD SELL
$ws1.Cells.Item(2, 1).Value = "2025-07-28 16:29"
$ws1.Cells.Item(2, 2).Value = "NZDUSD"
$ws1.Cells.Item(2, 3).Value = "SELL"
$ws1.Cells.Item(2, 4).Value = 0.59027
$ws1.Cells.Item(2, 5).Value = 0.59471
$ws1.Cells.Item(2, 6).Value = 0.5813
$ws1.Cells.Item(2, 7).Value = 0.05
$ws1.Cells.Item(2, 8).Value = "'" + "68.0%"
$ws1.Cells.Item(2, 9).Value = 2.02
$ws1.Cells.Item(2, 10).Value = "Active"

# Row 3: NZDUSD SELL
$ws1.Cells.Item(3, 1).Value = "2025-07-28 16:03"
$ws1.Cells.Item(3, 2).Value = "NZDUSD"
$ws1.Cells.Item(3, 3).Value = "SELL"
$ws1.Cells.Item(3, 4).Value = 0.59092
$ws1.Cells.Item(3, 5).Value = 0.59559
$ws1.Cells.Item(3, 6).Value = 0.58222
$ws1.Cells.Item(3, 7).Value = 0.01
$ws1.Cells.Item(3, 8).Value = "'" + "84.0%"
$ws1.Cells.Item(3, 9).Value = 1.86
$ws1.Cells.Item(3, 10).Value = "Active"

# Row 4: USDJPY SELL
$ws1.Cells.Item(4, 1).Value = "2025-07-28 15:53"
$ws1.Cells.Item(4, 2).Value = "USDJPY"
$ws1.Cells.Item(4, 3).Value = "SELL"
$ws1.Cells.Item(4, 4).Value = 148.7591
$ws1.Cells.Item(4, 5).Value = 148.97191
$ws1.Cells.Item(4, 6).Value = 148.00594
$ws1.Cells.Item(4, 7).Value = 0.08
$ws1.Cells.Item(4, 8).Value = "'" + "66.0%"
$ws1.Cells.Item(4, 9).Value = 3.54
$ws1.Cells.Item(4, 10).Value = "Active"

# Row 5: NZDUSD SELL
$ws1.Cells.Item(5, 1).Value = "2025-07-28 16:34"
$ws1.Cells.Item(5, 2).Value = "NZDUSD"
$ws1.Cells.Item(5, 3).Value = "SELL"
$ws1.Cells.Item(5, 4).Value = 0.59032
$ws1.Cells.Item(5, 5).Value = 0.5938
$ws1.Cells.Item(5, 6).Value = 0.58405
$ws1.Cells.Item(5, 7).Value = 0.02
$ws1.Cells.Item(5, 8).Value = "'" + "87.0%"
$ws1.Cells.Item(5, 9).Value = 1.8
$ws1.Cells.Item(5, 10).Value = "Active"

# Row 6: USDCHF SELL
$ws1.Cells.Item(6, 1).Value = "2025-07-28 16:07"
$ws1.Cells.Item(6, 2).Value = "USDCHF"
$ws1.Cells.Item(6, 3).Value = "SELL"
$ws1.Cells.Item(6, 4).Value = 0.88256
$ws1.Cells.Item(6, 5).Value = 0.88551
$ws1.Cells.Item(6, 6).Value = 0.87266
$ws1.Cells.Item(6, 7).Value = 0.09
$ws1.Cells.Item(6, 8).Value = "'" + "82.0%"
$ws1.Cells.Item(6, 9).Value = 3.36
$ws1.Cells.Item(6, 10).Value = "Active"

# Fix the Signal-column (C) fill: every row is now SELL, so copy the
# red SELL style (from C2) onto the rows that used to carry the green BUY style.
$ws1.Range("C2").Copy()
$ws1.Cells.Item(4, 3).PasteSpecial(-4122)
$ws1.Range("C2").Copy()
$ws1.Cells.Item(6, 3).PasteSpecial(-4122)

# Restore the plain bordered style (no quote-prefix) on the H-column cells
# after the apostrophe-forced text entry above.
$ws1.Range("A2").Copy()
$ws1.Cells.Item(2, 8).PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws1.Cells.Item(3, 8).PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws1.Cells.Item(4, 8).PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws1.Cells.Item(5, 8).PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws1.Cells.Item(6, 8).PasteSpecial(-4122)

# ============================================================
# Sheet 2: Summary Dashboard
# ============================================================
$ws2 = $wb.Worksheets.Item("Summary Dashboard")
$ws2.Range("B4").Value = 5
$ws2.Range("B5").Value = 6
$ws2.Range("B6").Value = 9
$ws2.Range("B7").Value = "77.6%"
$ws2.Range("B8").Value = "'" + "2.09"
$ws2.Range("B9").Value = "2025-07-28 16:19:05"

# Restore default (no explicit) style on B8 after the apostrophe-forced entry
$ws2.Range("B3").Copy()
$ws2.Range("B8").PasteSpecial(-4122)

# ============================================================
# Sheet 3: Signal History
# ============================================================
$ws3 = $wb.Worksheets.Item("Signal History")

# Row 2: AUDUSD BUY
$ws3.Cells.Item(2, 1).Value = "2025-07-28 16:38"
$ws3.Cells.Item(2, 2).Value = "AUDUSD"
$ws3.Cells.Item(2, 3).Value = "BUY"
$ws3.Cells.Item(2, 4).Value = 0.65446
$ws3.Cells.Item(2, 5).Value = 0.65779
$ws3.Cells.Item(2, 6).Value = 0.64612
$ws3.Cells.Item(2, 7).Value = 0.05
$ws3.Cells.Item(2, 8).Value = 0.75
$ws3.Cells.Item(2, 9).Value = 2.5
$ws3.Cells.Item(2, 10).Value = "Pending"

# Row 3: NZDUSD SELL
$ws3.Cells.Item(3, 1).Value = "2025-07-28 16:29"
$ws3.Cells.Item(3, 2).Value = "NZDUSD"
$ws3.Cells.Item(3, 3).Value = "SELL"
$ws3.Cells.Item(3, 4).Value = 0.59027
$ws3.Cells.Item(3, 5).Value = 0.59471
$ws3.Cells.Item(3, 6).Value = 0.5813
$ws3.Cells.Item(3, 7).Value = 0.05
$ws3.Cells.Item(3, 8).Value = 0.68
$ws3.Cells.Item(3, 9).Value = 2.02
$ws3.Cells.Item(3, 10).Value = "Active"

# Row 4: EURUSD BUY
$ws3.Cells.Item(4, 1).Value = "2025-07-28 16:31"
$ws3.Cells.Item(4, 2).Value = "EURUSD"
$ws3.Cells.Item(4, 3).Value = "BUY"
$ws3.Cells.Item(4, 4).Value = 1.10518
$ws3.Cells.Item(4, 5).Value = 1.10121
$ws3.Cells.Item(4, 6).Value = 1.11047
$ws3.Cells.Item(4, 7).Value = 0.03
$ws3.Cells.Item(4, 8).Value = 0.79
$ws3.Cells.Item(4, 9).Value = 1.33
$ws3.Cells.Item(4, 10).Value = "Pending"

# Row 5: NZDUSD SELL
$ws3.Cells.Item(5, 1).Value = "2025-07-28 16:03"
$ws3.Cells.Item(5, 2).Value = "NZDUSD"
$ws3.Cells.Item(5, 3).Value = "SELL"
$ws3.Cells.Item(5, 4).Value = 0.59092
$ws3.Cells.Item(5, 5).Value = 0.59559
$ws3.Cells.Item(5, 6).Value = 0.58222
$ws3.Cells.Item(5, 7).Value = 0.01
$ws3.Cells.Item(5, 8).Value = 0.84
$ws3.Cells.Item(5, 9).Value = 1.86
$ws3.Cells.Item(5, 10).Value = "Active"

# Row 6: USDJPY BUY
$ws3.Cells.Item(6, 1).Value = "2025-07-28 15:53"
$ws3.Cells.Item(6, 2).Value = "USDJPY"
$ws3.Cells.Item(6, 3).Value = "BUY"
$ws3.Cells.Item(6, 4).Value = 148.7591
$ws3.Cells.Item(6, 5).Value = 148.97191
$ws3.Cells.Item(6, 6).Value = 148.00594
$ws3.Cells.Item(6, 7).Value = 0.08
$ws3.Cells.Item(6, 8).Value = 0.66
$ws3.Cells.Item(6, 9).Value = 3.54
$ws3.Cells.Item(6, 10).Value = "Pending"

# Row 7: NZDUSD BUY
$ws3.Cells.Item(7, 1).Value = "2025-07-28 16:36"
$ws3.Cells.Item(7, 2).Value = "NZDUSD"
$ws3.Cells.Item(7, 3).Value = "BUY"
$ws3.Cells.Item(7, 4).Value = 0.58971
$ws3.Cells.Item(7, 5).Value = 0.58547
$ws3.Cells.Item(7, 6).Value = 0.5941
$ws3.Cells.Item(7, 7).Value = 0.05
$ws3.Cells.Item(7, 8).Value = 0.91
$ws3.Cells.Item(7, 9).Value = 1.04
$ws3.Cells.Item(7, 10).Value = "Filled"

# Row 8: NZDUSD SELL
$ws3.Cells.Item(8, 1).Value = "2025-07-28 16:34"
$ws3.Cells.Item(8, 2).Value = "NZDUSD"
$ws3.Cells.Item(8, 3).Value = "SELL"
$ws3.Cells.Item(8, 4).Value = 0.59032
$ws3.Cells.Item(8, 5).Value = 0.5938
$ws3.Cells.Item(8, 6).Value = 0.58405
$ws3.Cells.Item(8, 7).Value = 0.02
$ws3.Cells.Item(8, 8).Value = 0.87
$ws3.Cells.Item(8, 9).Value = 1.8
$ws3.Cells.Item(8, 10).Value = "Active"

# Row 9: AUDUSD SELL
$ws3.Cells.Item(9, 1).Value = "2025-07-28 16:19"
$ws3.Cells.Item(9, 2).Value = "AUDUSD"
$ws3.Cells.Item(9, 3).Value = "SELL"
$ws3.Cells.Item(9, 4).Value = 0.65982
$ws3.Cells.Item(9, 5).Value = 0.66341
$ws3.Cells.Item(9, 6).Value = 0.65043
$ws3.Cells.Item(9, 7).Value = 0.01
$ws3.Cells.Item(9, 8).Value = 0.92
$ws3.Cells.Item(9, 9).Value = 2.62
$ws3.Cells.Item(9, 10).Value = "Filled"

# Row 10: EURUSD BUY
$ws3.Cells.Item(10, 1).Value = "2025-07-28 16:34"
$ws3.Cells.Item(10, 2).Value = "EURUSD"
$ws3.Cells.Item(10, 3).Value = "BUY"
$ws3.Cells.Item(10, 4).Value = 1.10384
$ws3.Cells.Item(10, 5).Value = 1.09903
$ws3.Cells.Item(10, 6).Value = 1.10812
$ws3.Cells.Item(10, 7).Value = 0.08
$ws3.Cells.Item(10, 8).Value = 0.75
$ws3.Cells.Item(10, 9).Value = 0.89
$ws3.Cells.Item(10, 10).Value = "Pending"

# Row 11: NZDUSD BUY
$ws3.Cells.Item(11, 1).Value = "2025-07-28 15:54"
$ws3.Cells.Item(11, 2).Value = "NZDUSD"
$ws3.Cells.Item(11, 3).Value = "BUY"
$ws3.Cells.Item(11, 4).Value = 0.59185
$ws3.Cells.Item(11, 5).Value = 0.58911
$ws3.Cells.Item(11, 6).Value = 0.60011
$ws3.Cells.Item(11, 7).Value = 0.03
$ws3.Cells.Item(11, 8).Value = 0.65
$ws3.Cells.Item(11, 9).Value = 3.02
$ws3.Cells.Item(11, 10).Value = "Filled"

# Row 12: USDJPY BUY
$ws3.Cells.Item(12, 1).Value = "2025-07-28 16:11"
$ws3.Cells.Item(12, 2).Value = "USDJPY"
$ws3.Cells.Item(12, 3).Value = "BUY"
$ws3.Cells.Item(12, 4).Value = 150.20715
$ws3.Cells.Item(12, 5).Value = 149.90187
$ws3.Cells.Item(12, 6).Value = 150.6095
$ws3.Cells.Item(12, 7).Value = 0.07
$ws3.Cells.Item(12, 8).Value = 0.67
$ws3.Cells.Item(12, 9).Value = 1.32
$ws3.Cells.Item(12, 10).Value = "Filled"

# Row 13: NZDUSD SELL
$ws3.Cells.Item(13, 1).Value = "2025-07-28 16:06"
$ws3.Cells.Item(13, 2).Value = "NZDUSD"
$ws3.Cells.Item(13, 3).Value = "SELL"
$ws3.Cells.Item(13, 4).Value = 0.58912
$ws3.Cells.Item(13, 5).Value = 0.59309
$ws3.Cells.Item(13, 6).Value = 0.5834
$ws3.Cells.Item(13, 7).Value = 0.09
$ws3.Cells.Item(13, 8).Value = 0.66
$ws3.Cells.Item(13, 9).Value = 1.44
$ws3.Cells.Item(13, 10).Value = "Filled"

# Row 14: USDCHF SELL
$ws3.Cells.Item(14, 1).Value = "2025-07-28 16:07"
$ws3.Cells.Item(14, 2).Value = "USDCHF"
$ws3.Cells.Item(14, 3).Value = "SELL"
$ws3.Cells.Item(14, 4).Value = 0.88256
$ws3.Cells.Item(14, 5).Value = 0.88551
$ws3.Cells.Item(14, 6).Value = 0.87266
$ws3.Cells.Item(14, 7).Value = 0.09
$ws3.Cells.Item(14, 8).Value = 0.82
$ws3.Cells.Item(14, 9).Value = 3.36
$ws3.Cells.Item(14, 10).Value = "Active"

# Row 15: USDCAD BUY
$ws3.Cells.Item(15, 1).Value = "2025-07-28 16:29"
$ws3.Cells.Item(15, 2).Value = "USDCAD"
$ws3.Cells.Item(15, 3).Value = "BUY"
$ws3.Cells.Item(15, 4).Value = 1.36602
$ws3.Cells.Item(15, 5).Value = 1.36135
$ws3.Cells.Item(15, 6).Value = 1.37252
$ws3.Cells.Item(15, 7).Value = 0.02
$ws3.Cells.Item(15, 8).Value = 0.75
$ws3.Cells.Item(15, 9).Value = 1.39
$ws3.Cells.Item(15, 10).Value = "Pending"

# Row 16: NZDUSD SELL
$ws3.Cells.Item(16, 1).Value = "2025-07-28 16:44"
$ws3.Cells.Item(16, 2).Value = "NZDUSD"
$ws3.Cells.Item(16, 3).Value = "SELL"
$ws3.Cells.Item(16, 4).Value = 0.58751
$ws3.Cells.Item(16, 5).Value = 0.59064
$ws3.Cells.Item(16, 6).Value = 0.57763
$ws3.Cells.Item(16, 7).Value = 0.05
$ws3.Cells.Item(16, 8).Value = 0.92
$ws3.Cells.Item(16, 9).Value = 3.16
$ws3.Cells.Item(16, 10).Value = "Filled"

